$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Название компании"
$ws.Range("B1").Value = "Статус сертификата"
$ws.Range("C1").Value = "Номер сертификата"
$ws.Range("D1").Value = "ИНН"
$ws.Range("E1").Value = "Адрес"
$ws.Range("F1").Value = "Дата последнего аудита"
$ws.Range("G1").Value = "Стандарт"
$ws.Range("H1").Value = "Аккредитация"
$ws.Range("I1").Value = "Срок действия сертификата"

# New header cells G1/H1/I1 need the same left/center style as the rest of row 1.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:I1").PasteSpecial(-4122) | Out-Null

# --- Data row (row 2) ---
$ws.Range("A2").Value = '"ООО" Булка'
$ws.Range("A2").Style = "Обычный"
$ws.Range("B2").Value = "действующий"
$ws.Range("C2").Value = "03-00986"
$ws.Range("D2").Value = 474885885
$ws.Range("E2").Value = "СПб, ул.Вязовая, 14"
$ws.Range("F2").Value = 45407
$ws.Range("G2").Value = "ISO 9001"
$ws.Range("H2").Value = "UKAS (IAF)"
$ws.Range("I2").Value = 45412
$ws.Range("F2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I2").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'

# New data cells G2/H2 need the same left/center style as the rest of row 2.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("G2:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply values after the format paste (PasteSpecial(formats) shouldn't touch
# content, but this keeps things robust regardless of paste semantics).
$ws.Range("G1").Value = "Стандарт"
$ws.Range("H1").Value = "Аккредитация"
$ws.Range("I1").Value = "Срок действия сертификата"
$ws.Range("G2").Value = "ISO 9001"
$ws.Range("H2").Value = "UKAS (IAF)"

# --- Row 3: clear all old data (was a second certificate row) ---
$ws.Range("A3").Clear()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").Clear()
$ws.Range("E3").Clear()
$ws.Range("F3").ClearContents()

# --- Row 4: drop D4/E4 entirely, change F4's format to the date style ---
$ws.Range("D4").Clear()
$ws.Range("E4").Clear()
$ws.Range("F4").NumberFormat = "m/d/yy"

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 118
$ws.Range("D3").Select()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 27.0
$ws.Columns("B").ColumnWidth = 27.0
$ws.Columns("C").ColumnWidth = 22.0
$ws.Columns("D").ColumnWidth = 17.0
$ws.Columns("E").ColumnWidth = 21.166666666666668
$ws.Columns("F").ColumnWidth = 24.333333333333332
$ws.Columns("G").ColumnWidth = 14.333333333333334
$ws.Columns("H").ColumnWidth = 16.0
$ws.Columns("I").ColumnWidth = 27.333333333333332
